$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.582.18"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.689.03"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'314.13"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.3899"
$ws.Range("D8").Value = "'0.4033"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "'1.497"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "'52.83"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "'0.08743"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'25.22"
$ws.Range("E13").Value = "  +6.80%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'7.542"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").Value = "'7.953"
$ws.Range("E15").Value = "  -1.20%  "
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").Value = "1.698.56"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'98.62"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "'0.07106"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "'19.89"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "'7.279"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'14.26"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "24.571.37"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'3.013"
$ws.Range("E25").Value = "  -8.66%  "
$ws.Range("D26").Value = "'2.354"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'162.20"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'8.766"
$ws.Range("E29").Value = "  +14.58%  "
$ws.Range("D30").Value = "'137.08"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").Value = "'5.223"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").Value = "1.882.33"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'0.08826"
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("D34").Value = "'7.459"
$ws.Range("E34").Value = "  +4.73%  "
$ws.Range("D35").Value = "'1.039"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").Value = "'1.963"
$ws.Range("E36").Value = "  +4.23%  "
$ws.Range("D37").Value = "'0.02927"
$ws.Range("E37").Value = "  +7.48%  "
$ws.Range("D38").Value = "'0.2743"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "'10.79"
$ws.Range("E39").Value = "  -4.22%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'14.24"
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.09132"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'0.7853"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").Value = "'1.462"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'16.79"
$ws.Range("E44").Value = "  +4.61%  "
$ws.Range("D45").Value = "'0.7205"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "'2.595"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'4.202"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D49").Value = "'1.334"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").Value = "'137.94"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "'91.01"
$ws.Range("E51").Value = "  +0.57%  "
